$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(39).Insert()

$ws.Cells.Item(39, 1).Value = 11
$ws.Cells.Item(39, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(39, 3).Value = "Bíobío"
$ws.Cells.Item(39, 4).Value = 44624
$ws.Cells.Item(39, 5).Value = 8
$ws.Cells.Item(39, 6).Value = "Fruta"
$ws.Cells.Item(39, 7).Value = 100108
$ws.Cells.Item(39, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(39, 9).Value = 100108005
$ws.Cells.Item(39, 10).Value = "Piña"
$ws.Cells.Item(39, 11).Value = "Caramelo"
$ws.Cells.Item(39, 12).Value = "Segunda"
$ws.Cells.Item(39, 13).Value = 200
$ws.Cells.Item(39, 14).Value = 15000
$ws.Cells.Item(39, 15).Value = 16000
$ws.Cells.Item(39, 16).Value = 15600
$ws.Cells.Item(39, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(39, 18).Value = "Ecuador"
$ws.Cells.Item(39, 19).Value = 1114
$ws.Cells.Item(39, 20).Value = 14
